# Expenses Details.xlsx - add bill hyperlinks + CAN Analyzer / Battery Pack
# rows + Petrol Expenses / Person / Total summary columns, per commit:
# "Updated Expenses sheet for 1) Bill Links added 2) Bill for CAN Analyzer"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$missing = [System.Reflection.Missing]::Value

# ---------------------------------------------------------------------
# 1) New data rows 7 & 8 (Sr. No continues 6, 7)
# ---------------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("C7").Value = "CAN Analyzer"
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 3027
$ws.Range("F7").Value = "Parag"
$ws.Range("H7").Value = 0

$ws.Range("A8").Value = 7
$ws.Range("C8").Value = "Battery Pack and BMS"
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 18750

# ---------------------------------------------------------------------
# 2) Petrol Expenses column (H) for the existing 5 rows + new rows
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "Petrol Expenses"
$ws.Range("H2").Value = 100
$ws.Range("H3").Value = 100
$ws.Range("H4").Value = 100
$ws.Range("H5").Value = 100
$ws.Range("H6").Value = 100

# ---------------------------------------------------------------------
# 3) Totals block (rows 9-11)
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "Total"
$ws.Range("E9").Formula = "=SUM(E2:E8)"
$ws.Range("H9").Formula = "=SUM(H2:H8)"

$ws.Range("C10").Value = "Total Petrol + Expenses"
$ws.Range("E10").Formula = "=E9+H9"

$ws.Range("C11").Value = "Received"
$ws.Range("E11").Value = 5000
$ws.Range("G11").Value = "Balance"
$ws.Range("H11").Formula = "=E10-E11"

# ---------------------------------------------------------------------
# 4) Person / Total mini table (columns L & M)
# ---------------------------------------------------------------------
$ws.Range("M1").Value = "Total"
$ws.Range("L2").Value = "Akshay"
$ws.Range("M2").Formula = "=E2+E3+E4+E6"
$ws.Range("L3").Value = "Parag"
$ws.Range("M3").Formula = "=E5+E7"
$ws.Range("L1").Value = "Person"

# ---------------------------------------------------------------------
# 5) Bill hyperlinks (G2:G6) -> Bill_1 .. Bill_5
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G2"), "https://drive.google.com/drive/folders/Lithium-ion_battery_as_service/Bill_1", $missing, $missing) | Out-Null
$ws.Range("G2").Value = "Bill_1"

$ws.Hyperlinks.Add($ws.Range("G3"), "https://drive.google.com/drive/folders/Lithium-ion_battery_as_service/Bill_2", $missing, $missing) | Out-Null
$ws.Range("G3").Value = "Bill_2"

$ws.Hyperlinks.Add($ws.Range("G4"), "https://drive.google.com/drive/folders/Lithium-ion_battery_as_service/Bill_3", $missing, $missing) | Out-Null
$ws.Range("G4").Value = "Bill_3"

$ws.Hyperlinks.Add($ws.Range("G5"), "https://drive.google.com/drive/folders/Lithium-ion_battery_as_service/Bill_4", $missing, $missing) | Out-Null
$ws.Range("G5").Value = "Bill_4"

$ws.Hyperlinks.Add($ws.Range("G6"), "https://drive.google.com/drive/folders/Lithium-ion_battery_as_service/Bill_5", $missing, $missing) | Out-Null
$ws.Range("G6").Value = "Bill_5"

# carry the same (bordered + hyperlink-font) look into the two rows that
# don't have an actual bill attached, without creating a real hyperlink
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G7:G8").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 6) Formatting - reuse existing formats so we don't fork new styles
# ---------------------------------------------------------------------

# plain thin-bordered cells (same look as A2) -> new numeric/text cells
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7:A8").PasteSpecial(-4122) | Out-Null
$ws.Range("D7:D8").PasteSpecial(-4122) | Out-Null
$ws.Range("E7:E8").PasteSpecial(-4122) | Out-Null
$ws.Range("F7:F8").PasteSpecial(-4122) | Out-Null
$ws.Range("H2:H8").PasteSpecial(-4122) | Out-Null
$ws.Range("L2:L3").PasteSpecial(-4122) | Out-Null
$ws.Range("M2:M3").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:H9").PasteSpecial(-4122) | Out-Null
$ws.Range("A10:H10").PasteSpecial(-4122) | Out-Null
$ws.Range("A11:H11").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# date-style cells (same look as B2) -> blank but formatted B7:B8
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B7:B8").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# wrap-text detail cells (same look as C2) -> C7 / C8
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C7:C8").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# C10 label also wraps like the detail column
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# I8 / J8 - thin left border only next to the 18750 "funds" figure
$ws.Range("I8").Borders.Item(7).LineStyle = 1
$ws.Range("J8").Borders.Item(7).LineStyle = 0

# Header row - bold text on a yellow fill, like the rest of the header band
$headerRow = $ws.Range("A1:H1")
$headerRow.Font.Bold = $true
$headerRow.Interior.Color = 65535
$ws.Range("H1").Borders.Item(8).LineStyle = -4142

$ws.Range("L1:M1").Font.Bold = $true
$ws.Range("L1:M1").Interior.Color = 65535
$ws.Range("L1:M1").Borders.Item(7).LineStyle = 1
$ws.Range("L1:M1").Borders.Item(8).LineStyle = 1
$ws.Range("L1:M1").Borders.Item(9).LineStyle = 1
$ws.Range("L1:M1").Borders.Item(10).LineStyle = 1

# ---------------------------------------------------------------------
# 7) Column widths for the new columns
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 14.43
$ws.Columns.Item(13).ColumnWidth = 12.71

# ---------------------------------------------------------------------
# 8) Selection, matching where the author's cursor ended up
# ---------------------------------------------------------------------
$ws.Range("L7").Select() | Out-Null
